$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.084.03'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.92%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.377.76'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.53%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '555.22'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.17%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.88%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.631'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +1.97%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.370.02'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.56%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +5.36%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.97%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '53.55'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.75%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.27%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.14'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.49%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.916.67'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.71%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '18.31'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.61%  '
$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').NumberFormat = '@'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.118'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.14%  '
$ws.Range('B18').NumberFormat = '@'
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').NumberFormat = '@'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.378.57'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.55%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '64.955.77'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.84%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.20%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.03%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '457.58'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.02%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.19%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.32%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '14.16'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +5.62%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '87.68'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.86%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.20%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.73%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.23%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.03'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +3.05%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -2.07%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '63.20'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +7.33%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.44'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.31%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '576.36'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.71%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.70%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.62'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +3.41%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.44%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '35.61'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.13%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.36%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -2.74%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.100.53'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.05%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.16%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.75'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.68%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.19'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.39%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.05%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -3.84%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.00'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.05%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '140.40'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.68%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.53'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.09%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.13%  '
